$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45406 -> 45436)
$ws.Range("A1").Value = 45436

# Update the price list in column D (rows 33-41)
$ws.Range("D33").Value = 828
$ws.Range("D34").Value = 1037
$ws.Range("D35").Value = 1130
$ws.Range("D36").Value = 1408
$ws.Range("D37").Value = 1760
$ws.Range("D38").Value = 1446
$ws.Range("D39").Value = 1746
$ws.Range("D40").Value = 2153
$ws.Range("D41").Value = 2530
